$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D) updates ---
# Force text format on the whole Price column first so numeric-looking strings
# (e.g. "65.00", "236.10") are not silently coerced into floating point numbers,
# losing their original formatting. The style is reset back to Normal afterwards
# so the cells end up with no explicit style, matching the original workbook.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '39.751.81'
$ws.Range("D3").Value = '2.190.16'
$ws.Range("D5").Value = '290.87'
$ws.Range("D6").Value = '86.32'
$ws.Range("D10").Value = '30.01'
$ws.Range("D12").Value = '0.0775'
$ws.Range("D15").Value = '2.529.87'
$ws.Range("D16").Value = '13.65'
$ws.Range("D17").Value = '2.190.64'
$ws.Range("D18").Value = '0.722'
$ws.Range("D19").Value = '39.660.77'
$ws.Range("D20").Value = '0.0₃0878'
$ws.Range("D21").Value = '11.07'
$ws.Range("D22").Value = '5.69'
$ws.Range("D23").Value = '65.00'
$ws.Range("D24").Value = '236.10'
$ws.Range("D27").Value = '1.79'
$ws.Range("D28").Value = '22.99'
$ws.Range("D30").Value = '9.12'
$ws.Range("D31").Value = '155.88'
$ws.Range("D32").Value = '31.10'
$ws.Range("D34").Value = '4.88'
$ws.Range("D35").Value = '0.0702'
$ws.Range("D36").Value = '2.32'
$ws.Range("D39").Value = '0.0971'
$ws.Range("D40").Value = '1.66'
$ws.Range("D41").Value = '15.02'
$ws.Range("D42").Value = '2.117.24'
$ws.Range("D43").Value = '3.70'
$ws.Range("D45").Value = '0.0265'
$ws.Range("D46").Value = '9.68'
$ws.Range("D47").Value = '17.11'
$ws.Range("D49").Value = '2.403.23'
$ws.Range("D50").Value = '1.45'
$ws.Range("D51").Value = '87.62'

$priceRange.Style = "Normal"

# --- Other column updates (Coin name / Link / Volume%) ---
# These values are never numeric-looking (they contain letters, "%", "/", spaces)
# so Excel keeps them as plain text without any extra handling.
$ws.Range("E2").Value = '  -0.98%  '
$ws.Range("E3").Value = '  -2.37%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("E6").Value = '  -1.61%  '
$ws.Range("E7").Value = '  -1.80%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -2.75%  '
$ws.Range("E10").Value = '  -4.78%  '
$ws.Range("E11").Value = '  +5.97%  '
$ws.Range("E12").Value = '  -2.31%  '
$ws.Range("E13").Value = '  +2.34%  '
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("E15").Value = '  -1.86%  '
$ws.Range("E16").Value = '  -3.61%  '
$ws.Range("E17").Value = '  -1.50%  '
$ws.Range("E18").Value = '  -1.74%  '
$ws.Range("E19").Value = '  -0.94%  '
$ws.Range("E20").Value = '  -1.22%  '
$ws.Range("E21").Value = '  -1.81%  '
$ws.Range("E22").Value = '  -2.64%  '
$ws.Range("E23").Value = '  -1.00%  '
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("E26").Value = '  -2.14%  '
$ws.Range("E27").Value = '  -3.90%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("E29").Value = '  -8.12%  '
$ws.Range("E30").Value = '  -2.51%  '
$ws.Range("E31").Value = '  +2.18%  '
$ws.Range("E32").Value = '  -7.32%  '
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("E34").Value = '  -1.55%  '
$ws.Range("E35").Value = '  -3.15%  '
$ws.Range("E36").Value = '  -2.44%  '
$ws.Range("E37").Value = '  -1.28%  '
$ws.Range("E38").Value = '  -0.41%  '
$ws.Range("E39").Value = '  -3.48%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("E40").Value = '  -4.31%  '
$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("E41").Value = '  -9.33%  '
$ws.Range("E42").Value = '  +1.83%  '
$ws.Range("E43").Value = '  -4.15%  '
$ws.Range("E44").Value = '  -0.80%  '
$ws.Range("E45").Value = '  -1.84%  '
$ws.Range("E46").Value = '  -2.80%  '
$ws.Range("E47").Value = '  -7.33%  '
$ws.Range("E48").Value = '  +1.19%  '
$ws.Range("E49").Value = '  -1.40%  '
$ws.Range("E50").Value = '  -0.94%  '
$ws.Range("E51").Value = '  -2.19%  '

